# Auto-generated Excel COM-interop script applying numeric updates
# described by the upstream diff to Sheets/Cerberus_Profits.xlsx.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) receives a
# batch of cell-value writes; a handful of cells are cleared entirely
# (where the diff removed the <c> element) or newly populated (where
# the diff added one).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1857.3334  # H18
$ws.Cells.Item(18, 9).Value = 1840.2  # I18
$ws.Cells.Item(18, 11).Value = 1840.2  # K18
$ws.Cells.Item(18, 13).Value = -1556.2  # M18
$ws.Cells.Item(19, 8).Value = 1499.75  # H19
$ws.Cells.Item(19, 9).Value = 1750  # I19
$ws.Cells.Item(19, 10).Value = 1249.5  # J19
$ws.Cells.Item(19, 11).Value = 1750  # K19
$ws.Cells.Item(19, 12).Value = 1249.5  # L19
$ws.Cells.Item(19, 13).Value = -1575  # M19
$ws.Cells.Item(19, 14).Value = -1599.5  # N19
$ws.Cells.Item(38, 8).Value = 479.6154  # H38
$ws.Cells.Item(38, 9).Value = 479.6154  # I38
$ws.Cells.Item(38, 11).Value = 1438.8462  # K38
$ws.Cells.Item(38, 13).Value = -1066.8462  # M38
$ws.Cells.Item(57, 8).Value = 86446.28999999999  # H57
$ws.Cells.Item(57, 10).Value = 86446.28999999999  # J57
$ws.Cells.Item(57, 12).Value = 259338.87  # L57
$ws.Cells.Item(57, 14).Value = -260336.87  # N57
$ws.Cells.Item(103, 8).Value = 1361  # H103
$ws.Cells.Item(103, 9).Value = 1350  # I103
$ws.Cells.Item(103, 10).Value = 1368.3334  # J103
$ws.Cells.Item(103, 11).Value = 4050  # K103
$ws.Cells.Item(103, 12).Value = 4105.0002  # L103
$ws.Cells.Item(103, 13).Value = -3464  # M103
$ws.Cells.Item(103, 14).Value = -5277.0002  # N103
$ws.Cells.Item(132, 8).Value = 3222.84  # H132
$ws.Cells.Item(132, 9).Value = 2800.3  # I132
$ws.Cells.Item(132, 11).Value = 8400.900000000001  # K132
$ws.Cells.Item(132, 13).Value = -5870.900000000001  # M132
$ws.Cells.Item(137, 8).Value = 6279.2  # H137
$ws.Cells.Item(137, 9).Value = 2133  # I137
$ws.Cells.Item(137, 11).Value = 6399  # K137
$ws.Cells.Item(137, 13).Value = -3849  # M137
$ws.Cells.Item(141, 8).Value = 9151  # H141
$ws.Cells.Item(141, 9).Value = 5794.8  # I141
$ws.Cells.Item(141, 10).Value = 13945.571  # J141
$ws.Cells.Item(141, 11).Value = 17384.4  # K141
$ws.Cells.Item(141, 12).Value = 41836.713  # L141
$ws.Cells.Item(141, 13).Value = -12204.4  # M141
$ws.Cells.Item(141, 14).Value = -52196.713  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 805.39746  # H32
$ws.Cells.Item(32, 9).Value = 595.4776000000001  # I32
$ws.Cells.Item(32, 11).Value = 595.4776000000001  # K32
$ws.Cells.Item(32, 13).Value = -308.4776000000001  # M32
$ws.Cells.Item(61, 8).Value = 2121.2173  # H61
$ws.Cells.Item(61, 9).Value = 1620.4736  # I61
$ws.Cells.Item(61, 10).Value = 4499.75  # J61
$ws.Cells.Item(61, 11).Value = 1620.4736  # K61
$ws.Cells.Item(61, 12).Value = 4499.75  # L61
$ws.Cells.Item(61, 13).Value = -1408.4736  # M61
$ws.Cells.Item(61, 14).Value = -4923.75  # N61
$ws.Cells.Item(80, 8).Value = 99998.5  # H80
$ws.Cells.Item(80, 10).Value = 99998.5  # J80
$ws.Cells.Item(80, 12).Value = 99998.5  # L80
$ws.Cells.Item(80, 14).Value = -101994.5  # N80
$ws.Cells.Item(83, 8).Value = 99998.5  # H83
$ws.Cells.Item(83, 10).Value = 99998.5  # J83
$ws.Cells.Item(83, 12).Value = 299995.5  # L83
$ws.Cells.Item(83, 14).Value = -309979.5  # N83
$ws.Cells.Item(102, 8).Value = 3151.6956  # H102
$ws.Cells.Item(102, 9).Value = 3309.1428  # I102
$ws.Cells.Item(102, 11).Value = 3309.1428  # K102
$ws.Cells.Item(102, 13).Value = -1687.1428  # M102
$ws.Cells.Item(136, 8).Value = 2121.2173  # H136
$ws.Cells.Item(136, 9).Value = 1620.4736  # I136
$ws.Cells.Item(136, 10).Value = 4499.75  # J136
$ws.Cells.Item(136, 11).Value = 4861.4208  # K136
$ws.Cells.Item(136, 12).Value = 13499.25  # L136
$ws.Cells.Item(136, 13).Value = -2311.4208  # M136
$ws.Cells.Item(136, 14).Value = -18599.25  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(58, 8).Value = 49110.75  # H58
$ws.Cells.Item(58, 10).Value = 49110.75  # J58
$ws.Cells.Item(58, 12).Value = 49110.75  # L58
$ws.Cells.Item(58, 14).Value = -49698.75  # N58
$ws.Cells.Item(94, 8).Value = 6758.619  # H94
$ws.Cells.Item(94, 9).Value = 2036.3334  # I94
$ws.Cells.Item(94, 10).Value = 13055  # J94
$ws.Cells.Item(94, 11).Value = 2036.3334  # K94
$ws.Cells.Item(94, 12).Value = 13055  # L94
$ws.Cells.Item(94, 13).Value = -1585.3334  # M94
$ws.Cells.Item(94, 14).Value = -13957  # N94
$ws.Cells.Item(99, 8).Value = 1124  # H99
$ws.Cells.Item(99, 9).Value = 0  # I99
$ws.Cells.Item(99, 11).Value = 0  # K99
$ws.Cells.Item(99, 13).ClearContents()  # M99: was 1088
$ws.Cells.Item(105, 8).Value = 3016  # H105
$ws.Cells.Item(105, 9).Value = 2768.125  # I105
$ws.Cells.Item(105, 11).Value = 2768.125  # K105
$ws.Cells.Item(105, 13).Value = -1021.125  # M105
$ws.Cells.Item(107, 8).Value = 2969.6667  # H107
$ws.Cells.Item(107, 9).Value = 3090.875  # I107
$ws.Cells.Item(107, 10).Value = 2000  # J107
$ws.Cells.Item(107, 11).Value = 3090.875  # K107
$ws.Cells.Item(107, 12).Value = 2000  # L107
$ws.Cells.Item(107, 13).Value = -1170.875  # M107
$ws.Cells.Item(107, 14).Value = -5840  # N107
$ws.Cells.Item(134, 8).Value = 7558.517  # H134
$ws.Cells.Item(134, 9).Value = 7420.0586  # I134
$ws.Cells.Item(134, 10).Value = 7754.6665  # J134
$ws.Cells.Item(134, 11).Value = 22260.1758  # K134
$ws.Cells.Item(134, 12).Value = 23263.9995  # L134
$ws.Cells.Item(134, 13).Value = -19725.1758  # M134
$ws.Cells.Item(134, 14).Value = -28333.9995  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 410.3  # H22
$ws.Cells.Item(22, 9).Value = 362.16666  # I22
$ws.Cells.Item(22, 11).Value = 362.16666  # K22
$ws.Cells.Item(22, 13).Value = -12.16665999999998  # M22
$ws.Cells.Item(41, 8).Value = 32166.166  # H41
$ws.Cells.Item(41, 9).Value = 14000  # I41
$ws.Cells.Item(41, 10).Value = 35799.4  # J41
$ws.Cells.Item(41, 11).Value = 14000  # K41
$ws.Cells.Item(41, 12).Value = 35799.4  # L41
$ws.Cells.Item(41, 13).Value = -13572  # M41
$ws.Cells.Item(41, 14).Value = -36655.4  # N41
$ws.Cells.Item(50, 8).Value = 0  # H50
$ws.Cells.Item(50, 10).Value = 0  # J50
$ws.Cells.Item(50, 12).Value = 0  # L50
$ws.Cells.Item(50, 14).ClearContents()  # N50: was -61248
$ws.Cells.Item(51, 8).Value = 0  # H51
$ws.Cells.Item(51, 10).Value = 0  # J51
$ws.Cells.Item(51, 12).Value = 0  # L51
$ws.Cells.Item(51, 14).ClearContents()  # N51: was -61470
$ws.Cells.Item(53, 8).Value = 30000  # H53
$ws.Cells.Item(53, 10).Value = 30000  # J53
$ws.Cells.Item(53, 12).Value = 30000  # L53
$ws.Cells.Item(53, 14).Value = -31214  # N53
$ws.Cells.Item(58, 8).Value = 5721  # H58
$ws.Cells.Item(58, 9).Value = 1874.6  # I58
$ws.Cells.Item(58, 11).Value = 1874.6  # K58
$ws.Cells.Item(58, 13).Value = -1671.6  # M58
$ws.Cells.Item(59, 8).Value = 72749.75  # H59
$ws.Cells.Item(59, 9).Value = 46000  # I59
$ws.Cells.Item(59, 10).Value = 99499.5  # J59
$ws.Cells.Item(59, 11).Value = 46000  # K59
$ws.Cells.Item(59, 12).Value = 99499.5  # L59
$ws.Cells.Item(59, 13).Value = -44855  # M59
$ws.Cells.Item(59, 14).Value = -101789.5  # N59
$ws.Cells.Item(60, 8).Value = 25000  # H60
$ws.Cells.Item(60, 10).Value = 0  # J60
$ws.Cells.Item(60, 12).Value = 0  # L60
$ws.Cells.Item(60, 14).ClearContents()  # N60: was -56021
$ws.Cells.Item(61, 8).Value = 0  # H61
$ws.Cells.Item(61, 10).Value = 0  # J61
$ws.Cells.Item(61, 12).Value = 0  # L61
$ws.Cells.Item(61, 14).ClearContents()  # N61: was -60694
$ws.Cells.Item(136, 8).Value = 5721  # H136
$ws.Cells.Item(136, 9).Value = 1874.6  # I136
$ws.Cells.Item(136, 11).Value = 5623.799999999999  # K136
$ws.Cells.Item(136, 13).Value = -3073.799999999999  # M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(44, 8).Value = 1664  # H44
$ws.Cells.Item(44, 9).Value = 299.8  # I44
$ws.Cells.Item(44, 10).Value = 5074.5  # J44
$ws.Cells.Item(44, 11).Value = 899.4000000000001  # K44
$ws.Cells.Item(44, 12).Value = 15223.5  # L44
$ws.Cells.Item(44, 13).Value = -501.4000000000001  # M44
$ws.Cells.Item(44, 14).Value = -16019.5  # N44
$ws.Cells.Item(92, 8).Value = 394  # H92
$ws.Cells.Item(92, 10).Value = 394  # J92
$ws.Cells.Item(92, 12).Value = 1182  # L92
$ws.Cells.Item(92, 14).Value = -3678  # N92
$ws.Cells.Item(141, 8).Value = 31911.545  # H141
$ws.Cells.Item(141, 9).Value = 5783.778  # I141
$ws.Cells.Item(141, 11).Value = 17351.334  # K141
$ws.Cells.Item(141, 13).Value = -12171.334  # M141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 53747.375  # H123
$ws.Cells.Item(123, 10).Value = 53747.375  # J123
$ws.Cells.Item(123, 12).Value = 53747.375  # L123
$ws.Cells.Item(123, 14).Value = -58647.375  # N123

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1999  # H16
$ws.Cells.Item(16, 9).Value = 1999  # I16
$ws.Cells.Item(16, 11).Value = 1999  # K16
$ws.Cells.Item(16, 13).Value = -1829  # M16
$ws.Cells.Item(22, 8).Value = 982.3333  # H22
$ws.Cells.Item(22, 9).Value = 531.6667  # I22
$ws.Cells.Item(22, 10).Value = 1433  # J22
$ws.Cells.Item(22, 11).Value = 531.6667  # K22
$ws.Cells.Item(22, 12).Value = 1433  # L22
$ws.Cells.Item(22, 13).Value = -236.6667  # M22
$ws.Cells.Item(22, 14).Value = -2023  # N22
$ws.Cells.Item(27, 8).Value = 982.3333  # H27
$ws.Cells.Item(27, 9).Value = 531.6667  # I27
$ws.Cells.Item(27, 10).Value = 1433  # J27
$ws.Cells.Item(27, 11).Value = 531.6667  # K27
$ws.Cells.Item(27, 12).Value = 1433  # L27
$ws.Cells.Item(27, 13).Value = -424.6667  # M27
$ws.Cells.Item(27, 14).Value = -1647  # N27
$ws.Cells.Item(132, 8).Value = 2557.5  # H132
$ws.Cells.Item(132, 9).Value = 1598  # I132
$ws.Cells.Item(132, 10).Value = 4156.6665  # J132
$ws.Cells.Item(132, 11).Value = 4794  # K132
$ws.Cells.Item(132, 12).Value = 12469.9995  # L132
$ws.Cells.Item(132, 13).Value = -2264  # M132
$ws.Cells.Item(132, 14).Value = -17529.9995  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 134077.16  # H5
$ws.Cells.Item(5, 9).Value = 12375.375  # I5
$ws.Cells.Item(5, 11).Value = 12375.375  # K5
$ws.Cells.Item(5, 13).Value = -12263.375  # M5
$ws.Cells.Item(96, 8).Value = 2958.75  # H96
$ws.Cells.Item(96, 9).Value = 1324  # I96
$ws.Cells.Item(96, 10).Value = 3939.6  # J96
$ws.Cells.Item(96, 11).Value = 1324  # K96
$ws.Cells.Item(96, 12).Value = 3939.6  # L96
$ws.Cells.Item(96, 13).Value = 49  # M96
$ws.Cells.Item(96, 14).Value = -6685.6  # N96
$ws.Cells.Item(138, 8).Value = 89999.5  # H138
$ws.Cells.Item(138, 10).Value = 89999.5  # J138
$ws.Cells.Item(138, 12).Value = 89999.5  # L138
$ws.Cells.Item(138, 14).Value = -100279.5  # N138
